$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update grade entries - rows 2 and 3 (sheet rows 9 and 10) now earned full credit
$ws.Range("D9").Value = 1
$ws.Range("D10").Value = 1

# Update the active selection to reflect where the editor left off
$ws.Activate()
$ws.Range("G13").Select()
